# feat: uplift to the new version of pyxform/cht-conf (#27)
#
# Survey sheet: the "inputs" begin_group row no longer carries a
# NO_LABEL appearance value (column C), so clear it.
#
# Settings sheet: the form_id column is dropped entirely (it is no
# longer a recognised setting), so the whole column is deleted and
# version/style/namespaces shift one column to the left. The header
# comments are moved to match.

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# --- survey sheet -------------------------------------------------------
# Remove the stray "NO_LABEL" appearance value on the inputs begin_group row.
$survey.Range("C3").Clear()

# --- settings sheet ------------------------------------------------------
# Capture the header comments before the column shift moves the cells
# (but not the comments, which stay glued to their row/column position).
$titleComment = $settings.Range("A1").Comment().Text()
$versionComment = $settings.Range("C1").Comment().Text()
$pagesComment = $settings.Range("D1").Comment().Text()
$namespacesComment = $settings.Range("E1").Comment().Text()

# Drop the form_id column (B); version/style/namespaces shift left.
$settings.Columns.Item(2).Delete()

# Clear out whatever comments now sit in the header row (they did not
# travel with the cells during the column delete) and rebuild them in
# the correct, shifted positions.
foreach ($addr in @("A1", "B1", "C1", "D1", "E1")) {
    $existing = $settings.Range($addr).Comment()
    if ($existing -ne $null) {
        $existing.Delete()
    }
}

$settings.Range("A1").AddComment($titleComment)
$settings.Range("B1").AddComment($versionComment)
$settings.Range("C1").AddComment($pagesComment)
$settings.Range("D1").AddComment($namespacesComment)

# Refresh the cached version-stamp formula result.
$settings.Range("B2").Calculate()
$settings.Range("B1").Select() | Out-Null

# --- restore the survey sheet as the active tab/selection ---------------
$survey.Activate() | Out-Null
$survey.Range("C4").Select() | Out-Null
